$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) figures
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "70.226.74"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -2.58%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.533.58"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -3.85%  "

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "578.98"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -1.54%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "169.66"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -2.56%  "

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.519"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -0.13%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.533.55"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -3.81%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.162"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -5.31%  "

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -1.40%  "

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -2.00%  "

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.987.38"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -4.20%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "70.086.66"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -7.70%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "25.30"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -1.69%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.543.26"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -3.08%  "

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +0.70%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.37"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -6.25%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "353.36"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -5.56%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "3.94"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -3.15%  "

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.95%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "69.96"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "4.02"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -4.20%  "

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -2.48%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.657.87"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -4.18%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +1.37%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0915"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -3.85%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.90"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -4.57%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "467.64"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -4.14%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.76"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -1.93%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.11%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.119"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +3.20%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "156.93"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -2.69%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "19.01"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -3.58%  "

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +0.00%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "4.82"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -0.80%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.322"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -1.52%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.60"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -6.90%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.34"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -9.04%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.16"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -14.31%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "38.28"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -1.91%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "144.58"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -3.74%  "

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -1.10%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.52"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -3.51%  "

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -3.41%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0737"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -0.75%  "

